$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update homework scores (ДЗ_1..ДЗ_4, columns C:F) for several students.
# Row 4 - Баранов Алексей
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 2

# Row 6 - Брюн Феликс
$ws.Range("F6").Value = 2

# Row 7 - Герюгов Ислам
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 2

# Row 8 - Головин Антон
$ws.Range("F8").Value = 2

# Row 11 - Жуков Никита
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 2

# Row 14 - Корнилов Даниил
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 2
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 2

# Row 22 - Саитов Артур
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = 2

# Row 25 - Теплюк Дмитрий
$ws.Range("D25").Value = 2

# Row 27 - Французов Константин
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = 2
$ws.Range("F27").Value = 2

# Row 28 - Хабибулина Майя
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 2
$ws.Range("F28").Value = 2

# Move the view: scroll the frozen pane back up and select F1.
$ws.Range("F1").Select()
